$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the timetable data (rows 2-7) into rows 8-13 so the
# table fills an A4 page better when printed.
$data = @(
    @("10A", "25-10-2025", "Mathematics"),
    @("10A", "28-10-2025", "Science"),
    @("10A", "30-10-2025", "English"),
    @("10B", "25-10-2025", "Social Science"),
    @("10B", "28-10-2025", "Computer"),
    @("10B", "30-10-2025", "English")
)

$startRow = 8
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Update the selection to match the new active range.
$ws.Range("A8:C13").Select()
